$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 2063
$ws1.Range("F5").Value = 351
$ws1.Range("F6").Value = 600
$ws1.Range("F9").Value = 10597
$ws1.Range("F13").Value = 201
$ws1.Range("F14").Value = 414
$ws1.Range("F15").Value = 7483
$ws1.Range("F17").Value = 712
$ws1.Range("F18").Value = 233
$ws1.Range("F19").Value = 65
$ws1.Range("F20").Value = 3319

# Sheet "全部类型" (sheet4): update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2063
$ws4.Range("F5").Value = 351
$ws4.Range("F6").Value = 600
$ws4.Range("F12").Value = 10597
$ws4.Range("F16").Value = 201
$ws4.Range("F17").Value = 414
$ws4.Range("F18").Value = 7483
$ws4.Range("F20").Value = 712
$ws4.Range("F21").Value = 233
$ws4.Range("F22").Value = 65
$ws4.Range("F23").Value = 3319
